$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix name fields: first comma in these provider names should be a period ---
$nameCells = @("E20", "F20", "E49", "E50", "E60", "F60")
foreach ($addr in $nameCells) {
    $cell = $ws.Range($addr)
    $s = $cell.Text
    $idx = $s.IndexOf(",")
    if ($idx -ge 0) {
        $new = $s.Substring(0, $idx) + "." + $s.Substring($idx + 1)
        $origStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value = $new
        $cell.Style = $origStyle
    }
}

# --- Fix "Importe" column (H2:H130): convert "1.234,56" -> "1234.56" text ---
for ($r = 2; $r -le 130; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $s = $cell.Text
    $new = $s.Replace(".", "").Replace(",", ".")
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $new
    $cell.Style = $origStyle
}
